$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update metrics for row 3 (metrics_sim_with_priors.json)
$ws.Range("D3").Value = 0.8301886792452831
$ws.Range("E3").Value = 0.9845626072041166
$ws.Range("H3").Value = 0.5695676274944568
$ws.Range("I3").Value = 0.1131247965174388
$ws.Range("K3").Value = 1042.356775300171

$ws.Range("Q3").Value = 37
$ws.Range("R3").Value = 143
$ws.Range("S3").Value = 492
$ws.Range("T3").Value = 1080
$ws.Range("U3").Value = 1714
$ws.Range("V3").Value = 6596
$ws.Range("W3").Value = 6490
$ws.Range("X3").Value = 6141
$ws.Range("Y3").Value = 5553
$ws.Range("Z3").Value = 4919

$ws.Range("AF3").Value = 0.994422
$ws.Range("AG3").Value = 0.978441
$ws.Range("AH3").Value = 0.925825
$ws.Range("AI3").Value = 0.837178
$ws.Range("AJ3").Value = 0.741595
